$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 (everything from row 16 down shifts to 17 down).
$ws.Rows(16).Insert()

# Bring formats for the new row 16 in line with the rest of the table (copy from the
# row above for most columns, and from the row that just got pushed down to 17 for
# the currency column B, which uses a border variant unique to that row).
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("J15").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# New row 16: "Expansion Hub" / 289.62
$ws.Range("A16").Value = "Expansion Hub"
$ws.Range("B16").Value = 289.62

# Row 17 (previously row 16, still blank in A/B) now gets its own line item.
$ws.Range("A17").Value = "Extra Battery, Servo Adpater, Servo Extension"
$ws.Range("B17").Value = 61.3

$excel.CutCopyMode = 0
